$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# The template used Word fields ("{ m:userdoc 'zone1' }" / "{ m:enduserdoc }")
# to hold M2Doc tags. The parser was updated to a token-based rewriter
# (TokenIteratorFieldRewriterSplit) that expects those tags as plain literal
# text ("{m:userdoc 'zone1'}" / "{m:enduserdoc}") instead of real field
# codes, so the fldChar/instrText field plumbing is replaced with plain w:t
# runs carrying the same characters (the leading/trailing field-code spaces
# become the literal "{" / "}").

function Get-ParagraphOfPosition($doc, $pos) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if (($pos -ge $p.Range.Start) -and ($pos -lt $p.Range.End)) {
            return $p
        }
    }
    return $null
}

function Find-FieldByCode($doc, $pattern) {
    for ($i = 1; $i -le $doc.Fields.Count; $i++) {
        $candidate = $doc.Fields.Item($i)
        if ($candidate.Code.Text -match $pattern) {
            return $candidate
        }
    }
    return $null
}

# --- "{m:userdoc 'zone1'}" field -> literal text runs ---
$zoneField = Find-FieldByCode $d "userdoc"
$p2 = (Get-ParagraphOfPosition $d $zoneField.Code.Start).Range
$xml2 = '<w:p xmlns:w="' + $wNs + '" w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979">' `
    + '<w:r><w:t>{</w:t></w:r>' `
    + '<w:r><w:t>m</w:t></w:r>' `
    + '<w:r><w:t>:userdoc ''zone1''</w:t></w:r>' `
    + '<w:r><w:t xml:space="preserve">}</w:t></w:r>' `
    + '</w:p>'
$p2.InsertXML($xml2)

# --- "{m:enduserdoc}" field -> literal text runs (keep the _GoBack bookmark) ---
$endField = Find-FieldByCode $d "enduserdoc"
$p4 = (Get-ParagraphOfPosition $d $endField.Code.Start).Range
$xml4 = '<w:p xmlns:w="' + $wNs + '" w:rsidP="00F5495F" w:rsidR="00833091" w:rsidRDefault="00833091">' `
    + '<w:r><w:t>{m:</w:t></w:r>' `
    + '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' `
    + '<w:bookmarkEnd w:id="0"/>' `
    + '<w:r><w:t xml:space="preserve">enduserdoc}</w:t></w:r>' `
    + '</w:p>'
$p4.InsertXML($xml4)
